$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.467.93'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.573.01'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').Value = "'292.10"
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D8').Value = "'49.89"
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').Value = "'0.3414"
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').Value = "'0.07555"
$ws.Range('E11').Value = '  -1.24%  '
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').Value = "'6.051"
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').Value = "'6.975"
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = '1.579.82'
$ws.Range('E16').Value = '  +0.57%  '
$ws.Range('D17').Value = "'0.00001127"
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = "'91.30"
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('D19').Value = "'0.06763"
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('E21').Value = '  +1.29%  '
$ws.Range('D22').Value = "'16.36"
$ws.Range('E22').Value = '  -2.32%  '
$ws.Range('D23').Value = "'12.17"
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('D24').Value = '22.474.87'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').Value = "'2.663"
$ws.Range('E26').Value = '  +0.41%  '
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').Value = "'149.00"
$ws.Range('E28').Value = '  +1.31%  '
$ws.Range('D29').Value = "'5.059"
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('D30').Value = "'125.82"
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('D31').Value = '1.753.08'
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('D32').Value = "'1.080"
$ws.Range('E32').Value = '  +10.03%  '
$ws.Range('D33').Value = "'6.231"
$ws.Range('E33').Value = '  +1.11%  '
$ws.Range('D34').Value = "'2.019"
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('D35').Value = "'9.880"
$ws.Range('E35').Value = '  -2.40%  '
$ws.Range('D36').Value = "'0.08403"
$ws.Range('E36').Value = '  -0.87%  '
$ws.Range('D37').Value = "'0.02490"
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').Value = "'1.341"
$ws.Range('E39').Value = '  -2.69%  '
$ws.Range('D40').Value = "'0.06539"
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('D41').Value = "'5.470"
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('E42').Value = '  -0.66%  '
$ws.Range('D43').Value = "'0.6246"
$ws.Range('E43').Value = '  -2.19%  '
$ws.Range('D44').Value = "'14.07"
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D46').Value = "'3.812"
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = "'0.5837"
$ws.Range('E47').Value = '  -2.32%  '
$ws.Range('D48').Value = "'130.65"
$ws.Range('E48').Value = '  +4.73%  '
$ws.Range('D49').Value = "'2.083"
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('E50').Value = '  -5.17%  '
$ws.Range('D51').Value = "'0.07329"
$ws.Range('E51').Value = '  -0.03%  '

Write-Output "Applied cryptos update"
